$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '27.217.94'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.904.27'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  +0.29%  '
Set-TextValue $ws.Range("D5") '307.88'
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("E6").Value = '  +0.23%  '
Set-TextValue $ws.Range("D7") '0.5257'
$ws.Range("E7").Value = '  +0.50%  '
Set-TextValue $ws.Range("D8") '0.3829'
$ws.Range("E8").Value = '  +1.73%  '
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("E10").Value = '  +1.96%  '
$ws.Range("E11").Value = '  +0.40%  '
Set-TextValue $ws.Range("D12") '0.08079'
$ws.Range("E12").Value = '  -4.58%  '
Set-TextValue $ws.Range("D13") '95.95'
$ws.Range("E13").Value = '  +0.96%  '
Set-TextValue $ws.Range("D14") '5.368'
$ws.Range("E14").Value = '  +1.51%  '
$ws.Range("D15").Value = '1.784.04'
$ws.Range("E15").Value = '  -6.34%  '
$ws.Range("E16").Value = '  +0.23%  '
Set-TextValue $ws.Range("D17") '0.000008676'
$ws.Range("E17").Value = '  +0.48%  '
Set-TextValue $ws.Range("D18") '14.75'
$ws.Range("E18").Value = '  +1.37%  '
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").Value = '27.249.20'
Set-TextValue $ws.Range("D21") '5.120'
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("E22").Value = '  +2.36%  '
Set-TextValue $ws.Range("D23") '6.511'
$ws.Range("E23").Value = '  +1.31%  '
Set-TextValue $ws.Range("D24") '2.340'
$ws.Range("E24").Value = '  +2.63%  '
Set-TextValue $ws.Range("D25") '150.15'
$ws.Range("E25").Value = '  +1.89%  '
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("E27").Value = '  -0.56%  '
Set-TextValue $ws.Range("D28") '116.77'
$ws.Range("E28").Value = '  +1.69%  '
$ws.Range("E29").Value = '  +0.68%  '
Set-TextValue $ws.Range("D30") '4.877'
$ws.Range("E30").Value = '  -0.24%  '
Set-TextValue $ws.Range("D31") '0.09229'
$ws.Range("E31").Value = '  -0.32%  '
Set-TextValue $ws.Range("D32") '0.8177'
$ws.Range("E32").Value = '  +1.00%  '
Set-TextValue $ws.Range("D33") '0.05069'
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("E34").Value = '  -0.55%  '
Set-TextValue $ws.Range("D35") '2.986'
$ws.Range("E35").Value = '  +1.37%  '
Set-TextValue $ws.Range("D36") '3.362'
$ws.Range("E36").Value = '  -2.26%  '
Set-TextValue $ws.Range("D37") '2.725'
$ws.Range("E37").Value = '  +4.29%  '
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("E39").Value = '  +0.23%  '
Set-TextValue $ws.Range("D40") '1.085'
$ws.Range("E40").Value = '  +0.88%  '
Set-TextValue $ws.Range("D41") '9.018'
$ws.Range("E41").Value = '  -0.03%  '
Set-TextValue $ws.Range("D42") '6.626'
$ws.Range("E42").Value = '  -0.09%  '
Set-TextValue $ws.Range("D43") '117.03'
$ws.Range("E43").Value = '  +0.56%  '
Set-TextValue $ws.Range("D44") '0.1523'
$ws.Range("E44").Value = '  +0.67%  '
Set-TextValue $ws.Range("D45") '0.4931'
$ws.Range("E45").Value = '  +1.40%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D46") '1.002'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D47") '10.15'
$ws.Range("E47").Value = '  +0.39%  '
Set-TextValue $ws.Range("D48") '1.638'
$ws.Range("E48").Value = '  +1.43%  '
Set-TextValue $ws.Range("D49") '38.57'
$ws.Range("E49").Value = '  +2.90%  '
Set-TextValue $ws.Range("D50") '64.19'
$ws.Range("E50").Value = '  +0.29%  '
Set-TextValue $ws.Range("D51") '0.05968'
